$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1377.8889
$ws.Range("I32").Value = 1501
$ws.Range("J32").Value = 1342.7142
$ws.Range("K32").Value = 1501
$ws.Range("L32").Value = 1342.7142
$ws.Range("M32").Value = -1175
$ws.Range("N32").Value = -1994.7142

$ws.Range("H40").Value = 40001988
$ws.Range("I40").Value = 1450
$ws.Range("J40").Value = 47621140
$ws.Range("K40").Value = 1450
$ws.Range("L40").Value = 47621140
$ws.Range("M40").Value = -1275
$ws.Range("N40").Value = -47621490

$ws.Range("H43").Value = 2001400.6
$ws.Range("I43").Value = 2400.5
$ws.Range("J43").Value = 3334067.2
$ws.Range("K43").Value = 2400.5
$ws.Range("L43").Value = 3334067.2
$ws.Range("M43").Value = -2331.5
$ws.Range("N43").Value = -3334205.2

$ws.Range("H51").Value = 2263.4119
$ws.Range("I51").Value = 2320
$ws.Range("J51").Value = 2199.75
$ws.Range("K51").Value = 2320
$ws.Range("L51").Value = 2199.75
$ws.Range("M51").Value = -1836
$ws.Range("N51").Value = -3167.75

$ws.Range("H55").Value = 234.38461
$ws.Range("I55").Value = 195.18182
$ws.Range("K55").Value = 195.18182
$ws.Range("M55").Value = 18.81818000000001

$ws.Range("H137").Value = 10639916
$ws.Range("I137").Value = 1559.8064
$ws.Range("J137").Value = 31251732
$ws.Range("K137").Value = 4679.4192
$ws.Range("L137").Value = 93755196
$ws.Range("M137").Value = -2129.4192
$ws.Range("N137").Value = -93760296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 24473.25
$ws.Range("J113").Value = 24473.25
$ws.Range("L113").Value = 24473.25
$ws.Range("N113").Value = -33151.25

$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 1070.0869
$ws.Range("I122").Value = 1005.0909
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3015.2727
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227

$ws.Range("H94").Value = 1257.16
$ws.Range("I94").Value = 1262.6842
$ws.Range("J94").Value = 1239.6666
$ws.Range("K94").Value = 1262.6842
$ws.Range("L94").Value = 1239.6666
$ws.Range("M94").Value = -811.6841999999999
$ws.Range("N94").Value = -2141.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 293.5
$ws.Range("I22").Value = 258
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 258
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 92
$ws.Range("N22").Value = -1100

$ws.Range("H31").Value = 2085.0334
$ws.Range("I31").Value = 1006.1667
$ws.Range("J31").Value = 2354.75
$ws.Range("K31").Value = 1006.1667
$ws.Range("L31").Value = 2354.75
$ws.Range("M31").Value = -711.1667
$ws.Range("N31").Value = -2944.75

$ws.Range("H34").Value = 2085.0334
$ws.Range("I34").Value = 1006.1667
$ws.Range("J34").Value = 2354.75
$ws.Range("K34").Value = 1006.1667
$ws.Range("L34").Value = 2354.75
$ws.Range("M34").Value = -804.1667
$ws.Range("N34").Value = -2758.75

$ws.Range("M47").ClearContents()
$ws.Range("H47").Value = 9999
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 9999
$ws.Range("K47").Value = 0
$ws.Range("N47").Value = -11131

$ws.Range("H58").Value = 4788.967
$ws.Range("I58").Value = 918
$ws.Range("J58").Value = 8176.0625
$ws.Range("K58").Value = 918
$ws.Range("L58").Value = 8176.0625
$ws.Range("M58").Value = -715
$ws.Range("N58").Value = -8582.0625

$ws.Range("H105").Value = 2078.2273
$ws.Range("I105").Value = 1986.4286
$ws.Range("J105").Value = 2238.875
$ws.Range("K105").Value = 1986.4286
$ws.Range("L105").Value = 2238.875
$ws.Range("M105").Value = -239.4286
$ws.Range("N105").Value = -5732.875

$ws.Range("H136").Value = 4788.967
$ws.Range("I136").Value = 918
$ws.Range("J136").Value = 8176.0625
$ws.Range("K136").Value = 2754
$ws.Range("L136").Value = 24528.1875
$ws.Range("M136").Value = -204
$ws.Range("N136").Value = -29628.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1216.8586
$ws.Range("I68").Value = 630.63416
$ws.Range("J68").Value = 1631.2587
$ws.Range("K68").Value = 1891.90248
$ws.Range("L68").Value = 4893.7761
$ws.Range("M68").Value = -1080.90248
$ws.Range("N68").Value = -6515.7761

$ws.Range("H71").Value = 1216.8586
$ws.Range("I71").Value = 630.63416
$ws.Range("J71").Value = 1631.2587
$ws.Range("K71").Value = 5675.707439999999
$ws.Range("L71").Value = 14681.3283
$ws.Range("M71").Value = -1619.707439999999
$ws.Range("N71").Value = -22793.3283

$ws.Range("H107").Value = 126442.53
$ws.Range("I107").Value = 244.18518
$ws.Range("J107").Value = 223795.55
$ws.Range("K107").Value = 732.5555400000001
$ws.Range("L107").Value = 671386.6499999999
$ws.Range("M107").Value = 1187.44446
$ws.Range("N107").Value = -675226.6499999999

$ws.Range("H113").Value = 147590.55
$ws.Range("I113").Value = 483.33334
$ws.Range("J113").Value = 154380.11
$ws.Range("K113").Value = 1450.00002
$ws.Range("L113").Value = 463140.33
$ws.Range("M113").Value = 719.9999800000001
$ws.Range("N113").Value = -467480.33

$ws.Range("H131").Value = 11653425
$ws.Range("J131").Value = 2385.1143
$ws.Range("L131").Value = 7155.342900000001
$ws.Range("N131").Value = -17235.3429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2469.4375
$ws.Range("I97").Value = 2681.818
$ws.Range("K97").Value = 2681.818
$ws.Range("M97").Value = -2185.818

$ws.Range("H122").Value = 520944.44
$ws.Range("I122").Value = 650690.5600000001
$ws.Range("J122").Value = 1960
$ws.Range("K122").Value = 1952071.68
$ws.Range("L122").Value = 5880
$ws.Range("M122").Value = -1949621.68
$ws.Range("N122").Value = -10780

$ws.Range("H132").Value = 2204.182
$ws.Range("I132").Value = 1547.0322
$ws.Range("J132").Value = 3771.2307
$ws.Range("K132").Value = 4641.096600000001
$ws.Range("L132").Value = 11313.6921
$ws.Range("M132").Value = -2111.096600000001
$ws.Range("N132").Value = -16373.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 381.90475
$ws.Range("I46").Value = 324.16666
$ws.Range("J46").Value = 458.8889
$ws.Range("K46").Value = 324.16666
$ws.Range("L46").Value = 458.8889
$ws.Range("M46").Value = -136.16666
$ws.Range("N46").Value = -834.8888999999999

$ws.Range("H93").Value = 1145.0278
$ws.Range("J93").Value = 1775.4166
$ws.Range("L93").Value = 1775.4166
$ws.Range("N93").Value = -4271.4166

$ws.Range("H132").Value = 1245708.2
$ws.Range("I132").Value = 1738005.8
$ws.Range("J132").Value = 2009.3158
$ws.Range("K132").Value = 5214017.4
$ws.Range("L132").Value = 6027.9474
$ws.Range("M132").Value = -5211487.4
$ws.Range("N132").Value = -11087.9474

$ws.Range("H136").Value = 1636.5352
$ws.Range("I136").Value = 963.4091
$ws.Range("J136").Value = 2733.4814
$ws.Range("K136").Value = 2890.2273
$ws.Range("L136").Value = 8200.4442
$ws.Range("M136").Value = -340.2273
$ws.Range("N136").Value = -13300.4442

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1484.6608
$ws.Range("I132").Value = 1118.8695
$ws.Range("J132").Value = 3167.3
$ws.Range("K132").Value = 3356.6085
$ws.Range("L132").Value = 9501.900000000001
$ws.Range("M132").Value = -826.6085000000003
$ws.Range("N132").Value = -14561.9
